$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Cơ sở dữ liệu" moves up to row 3 (right after "Đánh giá và kiểm định"),
# pushing "Toán rời rạc" and "Lập trình Python" down one row each, and the
# class-count / total-student figures for those three modules are updated.

$ws.Range("A3").Value = "Cơ sở dữ liệu"
$ws.Range("B3").Value = 12

$ws.Range("A4").Value = "Toán rời rạc"
$ws.Range("B4").Value = 4

$ws.Range("A5").Value = "Lập trình Python"
$ws.Range("B5").Value = 8

# Column C ("Total Students") stores its numbers as text in this sheet, so
# force text entry (matching the existing C2 cell) rather than letting the
# numeric-looking strings get auto-converted to numbers.
$ws.Range("C3:C5").NumberFormat = "@"

$ws.Range("C3").Value = "430"
$ws.Range("C4").Value = "120"
$ws.Range("C5").Value = "275"

$ws.Range("C3:C5").Style = "Normal"
